$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39 ---
$ws.Range("A39").Value = 104203443
$ws.Range("B39").Value = 90139
$ws.Range("D39").Value = "CR"
$ws.Range("E39").Value = 1975
$ws.Range("F39").Value = "Liten kandelabersvamp"
$ws.Range("G39").Value = "Artomyces cristatus"
$ws.Range("H39").Value = "(Kauffman) Jülich"
$ws.Range("P39").Value = "Kroktjärnet O., Dls"
$ws.Range("Q39").Value = 318138.5606396351
$ws.Range("R39").Value = 6554822.488291118
$ws.Range("S39").Value = 10
$ws.Range("Y39").Value = "'2022-10-19"
$ws.Range("Y39").ClearFormats()
$ws.Range("Z39").Value = "00:00"
$ws.Range("AA39").Value = "'2022-10-19"
$ws.Range("AA39").ClearFormats()
$ws.Range("AB39").Value = "00:00"
$ws.Range("AI39").Value = "Barrskog"
$ws.Range("AO39").Value = "Låga"
$ws.Range("AW39").Value = "Tommy Solberg"
$ws.Range("AX39").Value = "Tommy Solberg"

# --- Row 40 ---
$ws.Range("A40").Value = 112387478
$ws.Range("B40").Value = 95020
$ws.Range("E40").Value = 2569
$ws.Range("F40").Value = "Stor revmossa"
$ws.Range("G40").Value = "Bazzania trilobata"
$ws.Range("H40").Value = "(L.) Gray"
$ws.Range("Q40").Value = 318275
$ws.Range("R40").Value = 6554953
$ws.Range("AJ40").ClearContents()
$ws.Range("AK40").ClearContents()
$ws.Range("AO40").ClearContents()

# --- Row 41 ---
$ws.Range("A41").Value = 112387492
$ws.Range("B41").Value = 90832
$ws.Range("D41").Value = "NT"
$ws.Range("E41").Value = 4368
$ws.Range("F41").Value = "Dofttaggsvamp"
$ws.Range("G41").Value = "Hydnellum suaveolens"
$ws.Range("H41").Value = "(Scop.:Fr.) P. Karst."
$ws.Range("P41").Value = "Stora Stickshöjden, Dls"
$ws.Range("Q41").Value = 318046
$ws.Range("R41").Value = 6554741
$ws.Range("S41").Value = 5
$ws.Range("Y41").Value = "'2023-09-26"
$ws.Range("Y41").ClearFormats()
$ws.Range("AA41").Value = "'2023-09-26"
$ws.Range("AA41").ClearFormats()
$ws.Range("AW41").Value = "Anton Larsson"
$ws.Range("AX41").Value = "Anton Larsson, Maria Johansson"
$ws.Range("J41").ClearContents()
$ws.Range("K41").ClearContents()
$ws.Range("N41").ClearContents()
$ws.Range("Z41").ClearContents()
$ws.Range("AB41").ClearContents()
$ws.Range("AF41").ClearContents()
$ws.Range("AI41").ClearContents()
$ws.Range("AO41").ClearContents()

# --- Row 42 ---
$ws.Range("A42").Value = 112387479
$ws.Range("B42").Value = 94340
$ws.Range("D42").Value = "LC"
$ws.Range("E42").Value = 2590
$ws.Range("F42").Value = "Kornknutmossa"
$ws.Range("G42").Value = "Odontoschisma denudatum"
$ws.Range("H42").Value = "(Mart.) Dumort"
$ws.Range("Q42").Value = 318301
$ws.Range("R42").Value = 6554977
$ws.Range("AJ42").Value = "tall"
$ws.Range("AK42").Value = "Pinus sylvestris"
$ws.Range("AO42").Value = "Pinus sylvestris"
